$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids numeric auto-conversion
# for values that look numeric, e.g. "676.23"), then strip the temporary text
# number-format override so the cell style stays the same as before the edit.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range("D2").Value = '69.603.69'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '3.696.14'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  -0.06%  '
Set-TextValue $ws.Range("D5") '676.23'
$ws.Range("E5").Value = '  -1.04%  '
Set-TextValue $ws.Range("D6") '160.95'
$ws.Range("E6").Value = '  +0.79%  '
$ws.Range("E8").Value = '  +0.64%  '
Set-TextValue $ws.Range("D10") '7.12'
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("E12").Value = '  +0.73%  '
Set-TextValue $ws.Range("D13") '32.56'
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").Value = '3.686.73'
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("D15").Value = '69.605.67'
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("E16").Value = '  +2.05%  '
Set-TextValue $ws.Range("D17") '16.01'
$ws.Range("E17").Value = '  +1.31%  '
Set-TextValue $ws.Range("D18") '6.47'
$ws.Range("E18").Value = '  +0.43%  '
Set-TextValue $ws.Range("D19") '470.06'
$ws.Range("E19").Value = '  +0.33%  '
Set-TextValue $ws.Range("D20") '9.82'
$ws.Range("E20").Value = '  -2.75%  '
Set-TextValue $ws.Range("D21") '0.651'
$ws.Range("E21").Value = '  +0.79%  '
Set-TextValue $ws.Range("D22") '80.62'
$ws.Range("E22").Value = '  +1.28%  '
$ws.Range("D23").Value = '3.845.70'
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("E25").Value = '  +3.69%  '
Set-TextValue $ws.Range("D26") '10.88'
$ws.Range("E26").Value = '  -0.32%  '
Set-TextValue $ws.Range("D27") '9.12'
$ws.Range("E27").Value = '  -0.34%  '
Set-TextValue $ws.Range("D28") '2.70'
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("E29").Value = '  +1.70%  '
$ws.Range("E30").Value = '  +0.42%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D31") '1.01'
$ws.Range("E31").Value = '  +0.68%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D32") '6.58'
$ws.Range("E32").Value = '  +0.19%  '
Set-TextValue $ws.Range("D33") '26.96'
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("D34").Value = '3.690.22'
$ws.Range("E34").Value = '  +0.84%  '
Set-TextValue $ws.Range("D35") '0.163'
$ws.Range("E35").Value = '  +0.18%  '
Set-TextValue $ws.Range("D36") '8.46'
$ws.Range("E36").Value = '  +3.98%  '
Set-TextValue $ws.Range("D37") '6.22'
$ws.Range("E37").Value = '  +1.98%  '
$ws.Range("E38").Value = '  +0.00%  '
Set-TextValue $ws.Range("D39") '2.24'
$ws.Range("E39").Value = '  -1.95%  '
$ws.Range("E40").Value = '  -0.04%  '
Set-TextValue $ws.Range("D41") '0.0901'
$ws.Range("E41").Value = '  +0.30%  '
Set-TextValue $ws.Range("D42") '167.32'
$ws.Range("E42").Value = '  +0.76%  '
Set-TextValue $ws.Range("D43") '0.944'
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("E44").Value = '  -1.19%  '
$ws.Range("E45").Value = '  +2.24%  '
Set-TextValue $ws.Range("D46") '28.09'
$ws.Range("E46").Value = '  -0.60%  '
Set-TextValue $ws.Range("D47") '0.000278'
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("E49").Value = '  -3.22%  '
Set-TextValue $ws.Range("D50") '7.88'
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("E51").Value = '  +1.85%  '
